# Update the cached "datetimeFigureOut" footer field (Date Placeholder)
# from 16/04/2022 to 12/05/2022 on the slide master and on every slide
# layout that carries it.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($container, $newDate) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$newDate = "12/05/2022"

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master $newDate

# Every slide layout hanging off the master (accessed via the master's
# CustomLayouts collection so each index resolves to its own layout).
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout $newDate
}
